$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 13 (shifts rows 13-23 down to 14-24),
# carrying row heights and formatting down with them.
$ws.Rows.Item(13).Insert()

# The inserted row picks up a stray styled-but-empty cell in column A;
# remove it completely so row 13 has no A cell, matching the target layout.
$ws.Range("A13").Clear()

# Row 10
$ws.Range("B10").Value = @'
Desenvolver um projeto sobre tema de Engenharia de Produção, similar a situações que os alunos irão encontrar na vida real no efetivo exercício de sua profissão, 
Aplicar e integrar conhecimentos adquiridos em demais disciplinas de seu curso
Desenvolver competências técnicas, as relacionadas ao projeto em si, bem como competências transversais (habilidades e atitudes), num ambiente de aprendizagem baseado em PBL (Project-Baed Learning).
'@
$ws.Range("C10").Value = @'
Desenvolver um projeto sobre tema de Engenharia de Produção, similar a situações que os alunos irão encontrar na vida real no efetivo exercício de sua profissão, 
Aplicar e integrar conhecimentos adquiridos em demais disciplinas de seu curso
Desenvolver competências técnicas, as relacionadas ao projeto em si, bem como competências transversais (habilidades e atitudes), num ambiente de aprendizagem baseado em PBL (Project-Baed Learning).
'@

# Row 13
$ws.Range("B13").Value = '5840560 - Marco Antonio Carvalho Pereira'
$ws.Range("C13").Value = '5840560 - Marco Antonio Carvalho Pereira'

# Row 14
$ws.Range("B14").Value = 'Tópicos que abordem o tema do projeto de seu planejamento a execução.'
$ws.Range("C14").Value = 'Tópicos que abordem o tema do projeto de seu planejamento a execução.'

# Row 16
$ws.Range("B16").Value = @'
Noções de Gestão de Projetos
Organização do tempo: dimensão pessoal;
Técnicas para a realização de apresentações;
Noções de Aprendizagem Baseada em Projetos
Trabalho em Grupo, Equipes e times. 
Postura e Ética Profissional
Técnicas para redação de relatório técnico;
Tutoria de projetos.
Assuntos Técnicos específicos relacionados com o tema do projeto.
'@
$ws.Range("C16").Value = @'
Noções de Gestão de Projetos
Organização do tempo: dimensão pessoal;
Técnicas para a realização de apresentações;
Noções de Aprendizagem Baseada em Projetos
Trabalho em Grupo, Equipes e times. 
Postura e Ética Profissional
Técnicas para redação de relatório técnico;
Tutoria de projetos.
Assuntos Técnicos específicos relacionados com o tema do projeto.
'@

# Row 19
$ws.Range("B19").Value = @'
O método utilizado tem por fundamento a Aprendizagem Baseada em Projetos (PBL) que visa desenvolver as competências técnicas relativas ao tema do projeto, bem como competências transversais, tais como: aprender a aprender, trabalho em equipe, relacionamento interpessoal, aspectos de liderança e capacidade de comunicação, dentre outras.

Os alunos serão divididos em grupos que desenvolverão um projeto durante o semestre relacionado a um tema de Engenharia de Produção, similar ao que eles irão encontrar na vida real no efetivo exercício de sua profissão. 
Cada grupo deverá buscar o conhecimento prático necessário para ser aplicado no desenvolvimento do projeto.
As aulas ocorrerão: 1) através de uma reunião da equipe de trabalho para tratar do projeto, e  2) palestras e dinâmicas relativas ao tema do projeto, conduzidas por professores  ou profissionais de empresas.
'@
$ws.Range("C19").Value = @'
O método utilizado tem por fundamento a Aprendizagem Baseada em Projetos (PBL) que visa desenvolver as competências técnicas relativas ao tema do projeto, bem como competências transversais, tais como: aprender a aprender, trabalho em equipe, relacionamento interpessoal, aspectos de liderança e capacidade de comunicação, dentre outras.

Os alunos serão divididos em grupos que desenvolverão um projeto durante o semestre relacionado a um tema de Engenharia de Produção, similar ao que eles irão encontrar na vida real no efetivo exercício de sua profissão. 
Cada grupo deverá buscar o conhecimento prático necessário para ser aplicado no desenvolvimento do projeto.
As aulas ocorrerão: 1) através de uma reunião da equipe de trabalho para tratar do projeto, e  2) palestras e dinâmicas relativas ao tema do projeto, conduzidas por professores  ou profissionais de empresas.
'@

# Row 20
$ws.Range("B20").Value = @'
A nota será individual e será a média ponderada de componentes do projeto, tais como: Projeto Preliminar, Projeto Final, envolvimento do aluno com o projeto, Avaliação dos Pares, Apresentação de Trabalhos, dentre outros.
O detalhamento dos pesos para ponderação da média da disciplina será definido por uma equipe de professores que atuarão na coordenação da disciplina.
'@
$ws.Range("C20").Value = @'
A nota será individual e será a média ponderada de componentes do projeto, tais como: Projeto Preliminar, Projeto Final, envolvimento do aluno com o projeto, Avaliação dos Pares, Apresentação de Trabalhos, dentre outros.
O detalhamento dos pesos para ponderação da média da disciplina será definido por uma equipe de professores que atuarão na coordenação da disciplina.
'@

# Row 21
$ws.Range("B21").Value = 'Não há recuperação'
$ws.Range("C21").Value = 'Não há recuperação'

# Row 22
$ws.Range("B22").Value = @'
Artigos sobre metodologias ativas de aprendizagem e  Project Based Learning.
Livros e Artigos científicos relacionados com o tema do projeto.
'@
$ws.Range("C22").Value = @'
Artigos sobre metodologias ativas de aprendizagem e  Project Based Learning.
Livros e Artigos científicos relacionados com o tema do projeto.
'@
